$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update status from "Panding" to "Kualifikasi" for the rows that were
# qualified (K10, K12, K14, K15, K18).
$rows = @(10, 12, 14, 15, 18)
foreach ($r in $rows) {
    $ws.Range("K$r").Value = "Kualifikasi"
}

# Column K widened slightly to fit the new text and marked as best-fit.
$ws.Columns.Item(11).ColumnWidth = 13.6328125

# Update the view: scroll so column C is the left-most visible column and
# select G8 instead of C30.
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("G8").Select()
